$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.090.01"
$ws.Range("E2").Value = "  -0.45%  "

# Row 3
$ws.Range("D3").Value = "1.784.35"
$ws.Range("E3").Value = "  -2.40%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.71"
$ws.Range("E5").Value = "  +0.31%  "

# Row 6
$ws.Range("E6").Value = "  -1.43%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.72"
$ws.Range("E8").Value = "  +2.39%  "

# Row 9
$ws.Range("E9").Value = "  -2.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0712"
$ws.Range("E10").Value = "  -1.03%  "

# Row 11
$ws.Range("E11").Value = "  +0.60%  "

# Row 12
$ws.Range("D12").Value = "2.041.65"
$ws.Range("E12").Value = "  -2.60%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.98"
$ws.Range("E13").Value = "  +1.71%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.775.96"
$ws.Range("E14").Value = "  -3.10%  "

# Row 15
$ws.Range("D15").Value = "34.025.23"
$ws.Range("E15").Value = "  -0.81%  "

# Row 16
$ws.Range("E16").Value = "  -3.47%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.15"
$ws.Range("E17").Value = "  -4.22%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.78"
$ws.Range("E18").Value = "  -2.75%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.02"
$ws.Range("E19").Value = "  -2.67%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0788"
$ws.Range("E20").Value = "  -0.45%  "

# Row 21
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.77"
$ws.Range("E22").Value = "  -3.53%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  -4.24%  "

# Row 24
$ws.Range("E24").Value = "  -2.97%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.17"
$ws.Range("E25").Value = "  -0.18%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.33"
$ws.Range("E26").Value = "  -1.94%  "

# Row 27
$ws.Range("E27").Value = "  -2.96%  "

# Row 28
$ws.Range("E28").Value = "  -2.34%  "

# Row 29
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("E30").Value = "  +0.58%  "

# Row 31
$ws.Range("E31").Value = "  -4.64%  "

# Row 32
$ws.Range("E32").Value = "  -4.12%  "

# Row 33
$ws.Range("E33").Value = "  -1.86%  "

# Row 34
$ws.Range("E34").Value = "  -4.88%  "

# Row 35
$ws.Range("D35").Value = "1.391.62"
$ws.Range("E35").Value = "  -3.71%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.646"
$ws.Range("E36").Value = "  +0.06%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.05"
$ws.Range("E37").Value = "  -1.67%  "

# Row 38
$ws.Range("E38").Value = "  -2.21%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("E39").Value = "  +3.02%  "

# Row 40
$ws.Range("E40").Value = "  -0.26%  "

# Row 41
$ws.Range("E41").Value = "  -4.89%  "

# Row 42
$ws.Range("E42").Value = "  -2.67%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "78.12"
$ws.Range("E43").Value = "  -4.40%  "

# Row 44
$ws.Range("D44").Value = "0.0₆0143"
$ws.Range("E44").Value = "  +14.61%  "

# Row 45
$ws.Range("E45").Value = "  +2.68%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "108.28"
$ws.Range("E46").Value = "  +1.84%  "

# Row 47
$ws.Range("E47").Value = "  -0.44%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.47"
$ws.Range("E48").Value = "  +5.07%  "

# Row 49
$ws.Range("E49").Value = "  -4.28%  "

# Row 50
$ws.Range("D50").Value = "1.941.19"
$ws.Range("E50").Value = "  -2.69%  "

# Row 51
$ws.Range("E51").Value = "  +0.02%  "

